$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B labels shift because two new entries ("Holden", "Rizzie Spiral")
#     were spliced into the lookup list used to build this column, and
#     "Thomas Hex" was renamed to "Matthies Hex" ---
$ws.Range("B4").Value  = "Holden"
$ws.Range("B5").Value  = "Rizzie Spiral"
$ws.Range("B6").Value  = "RotRing OmegaMax-90"
$ws.Range("B7").Value  = "Equal Angle"
$ws.Range("B8").Value  = "Tilt Rotate"
$ws.Range("B9").Value  = "CLR"
$ws.Range("B10").Value = "Rizzie Hex"
$ws.Range("B11").Value = "Matthies Hex"
$ws.Range("B12").Value = "Tilt Rotate_Partial"
$ws.Range("B13").Value = "RotRing OmegaMax-60"
$ws.Range("B14").Value = "Equal Angle_Partial"
$ws.Range("B15").Value = "Rizzie Hex_Partial"
$ws.Range("B16").Value = "ND Single"
$ws.Range("B17").Value = "RD Single"
$ws.Range("B18").Value = "TD Single"
$ws.Range("B19").Value = "Morris Single"
$ws.Range("B20").Value = "Ring Perpendicular to ND"
$ws.Range("B21").Value = "Ring Perpendicular to RD"
$ws.Range("B22").Value = "Ring Perpendicular to TD"
$ws.Range("B23").Value = "OffsetFTD"
$ws.Range("B24").Value = "OffsetATD"
$ws.Range("B25").Value = "OffsetF45"
$ws.Range("B26").Value = "OffsetA45"
$ws.Range("B27").Value = "OffsetFRD"
$ws.Range("B28").Value = "OffsetARD"
$ws.Range("B29").Value = "Gaussian Quadrature"

# --- Append two new data rows (30 and 31), continuing the existing pattern ---
# Row 30: index 28 -> "Michael-CCHex"
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "Michael-CCHex"
$ws.Range("C30:T30").Value = 1

# Row 31: index 29 -> "Michael-SNHex"
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "Michael-SNHex"
$ws.Range("C31:T31").Value = 1

# Match the formatting of column A used by the rest of the table
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
